# Adding Code For Payer Insurance Creation
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pepID")

$values = @(
    "PEP_ID-2009581",
    "PEP_ID-2009582",
    "PEP_ID-2009583",
    "PEP_ID-2009584",
    "PEP_ID-2009589",
    "PEP_ID-2009594",
    "PEP_ID-2009599",
    "PEP_ID-2009602",
    "PEP_ID-2009603",
    "PEP_ID-2009604",
    "PEP_ID-2009605",
    "PEP_ID-2009606",
    "PEP_ID-2009607",
    "PEP_ID-2009609",
    "PEP_ID-2009611",
    "PEP_ID-2009612",
    "PEP_ID-2009630",
    "PEP_ID-2009631",
    "PEP_ID-2009632",
    "PEP_ID-2009633",
    "PEP_ID-2009642",
    "PEP_ID-2009643",
    "PEP_ID-2009647",
    "PEP_ID-2009648",
    "PEP_ID-2009649",
    "PEP_ID-2009650",
    "PEP_ID-2009651"
)

$row = 2
foreach ($val in $values) {
    $ws.Cells.Item($row, 1).Value = $val
    $row = $row + 1
}
